$d = $word.ActiveDocument

# 1) Remove the paragraph "आकान, आत्मा के वरदान..." (italic term list after the "आ" heading)
$d.Paragraphs.Item(11).Range.Delete()

# 2) Remove the paragraph "This PDF version is provided under the same license."
$d.Paragraphs.Item(6).Range.Delete()

# 3) Rewrite the license/attribution paragraph (previously bold key-term + hyperlinked
#    license text) with the new Biblica Study Notes attribution text.
$p5 = $d.Paragraphs.Item(5)
$pkg5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr/><w:r><w:rPr><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr></w:r><w:r><w:rPr><w:b/><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr><w:t>Biblica Study Notes (Key Terms)</w:t></w:r><w:r><w:rPr><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr><w:t xml:space="preserve"> © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. </w:t></w:r><w:r><w:rPr><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr><w:t>Biblica Study Notes</w:t></w:r><w:r><w:rPr><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr><w:t xml:space="preserve"> has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.</w:t></w:r><w:r><w:rPr><w:lang w:val="hi_IN" w:bidi="hi_IN"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5.Range.InsertXML($pkg5)

# 4) Remove the "License Information" heading paragraph
$d.Paragraphs.Item(4).Range.Delete()
